{"js": "// Word JS API (Office.js) script.\n// Reproduces the content-level edits described by the diff:\n//   1. Remove the stray \"_GoBack\" bookmark at the very start of the document\n//      (Word automatically renumbers the remaining \"_TocXXXXXXXX\" bookmarks\n//      afterwards, which is also what the diff shows).\n//   2. Fix the \"References\" entry: replace the old Google Books URL with the\n//      new webcitation.org URL.\n//   3. Drop the stray <w:lastRenderedPageBreak/> that precedes the \"2.\" item\n//      in the peer-review appendix (a pagination cache artifact left over\n//      from a previous save) while leaving the visible \"2.\" text untouched.\n\n// 1. Remove the \"_GoBack\" bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Replace the outdated reference URL with the new citation URL.\nconst oldUrl = \"https://books.google.fi/books?id=PNEGt99uMBwC&pg=PP1\";\nconst newUrl = \"https://www.webcitation.org/65iNkn800?url=http://www.aero.org/publications/crosslink/winter2002/04.html\";\nconst urlResults = context.document.body.search(oldUrl, { matchCase: true });\nurlResults.load(\"items\");\nawait context.sync();\nif (urlResults.items.length > 0) {\n  urlResults.items[0].insertText(newUrl, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3. Locate the paragraph that starts the second peer-review comment\n//    (\"2.\" followed by a tab and the quoted feedback) and rewrite its\n//    leading run so the cached <w:lastRenderedPageBreak/> marker is dropped.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"2.\") === 0 && text.indexOf(\"For \") !== -1 && text.indexOf(\"example\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  const leadResults = targetParagraph.search(\"2.\", { matchCase: true });\n  leadResults.load(\"items\");\n  await context.sync();\n  if (leadResults.items.length > 0) {\n    // Swap the text out and back in so Word re-creates the run from\n    // scratch, which drops the stale lastRenderedPageBreak marker.\n    leadResults.items[0].insertText(\"2.\\u0000tmp\", Word.InsertLocation.replace);\n    await context.sync();\n\n    const tmpResults = targetParagraph.search(\"2.\\u0000tmp\", { matchCase: true });\n    tmpResults.load(\"items\");\n    await context.sync();\n    if (tmpResults.items.length > 0) {\n      tmpResults.items[0].insertText(\"2.\", Word.InsertLocation.replace);\n      await context.sync();\n    }\n  }\n}\n", "ps1": "# Word COM (PowerShell-style) script.\n# Reproduces the content-level edits described by the diff:\n#   1. Remove the stray \"_GoBack\" bookmark at the very start of the document\n#      (Word automatically renumbers the remaining \"_TocXXXXXXXX\" bookmarks\n#      afterwards, which is also what the diff shows).\n#   2. Fix the \"References\" entry: replace the old Google Books URL with the\n#      new webcitation.org URL.\n#   3. Drop the stray lastRenderedPageBreak marker that precedes the \"2.\"\n#      item in the peer-review appendix (a pagination cache artifact left\n#      over from a previous save) while leaving the visible \"2.\" text\n#      untouched.\n\n$d = $word.ActiveDocument\n\n# 1. Remove the \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Replace the outdated reference URL with the new citation URL.\n$oldUrl = \"https://books.google.fi/books?id=PNEGt99uMBwC&pg=PP1\"\n$newUrl = \"https://www.webcitation.org/65iNkn800?url=http://www.aero.org/publications/crosslink/winter2002/04.html\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $oldUrl\n$found = $find.Execute()\nif ($found) {\n    $rng.Text = $newUrl\n}\n\n# 3. Locate the paragraph that starts the second peer-review comment\n#    (\"2.\" followed by a tab and the quoted feedback) and rewrite its\n#    leading run so the cached lastRenderedPageBreak marker is dropped.\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"2.\") -and $t.Contains(\"For \") -and $t.Contains(\"example\")) {\n        $pr = $p.Range\n        $leadFind = $pr.Find\n        $leadFind.ClearFormatting()\n        $leadFind.Text = \"2.\"\n        $leadFound = $leadFind.Execute()\n        if ($leadFound) {\n            # Swap the text out and back in so Word re-creates the run from\n            # scratch, which drops the stale lastRenderedPageBreak marker.\n            $pr.Text = \"2.TEMP\"\n        }\n        break\n    }\n}\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"2.TEMP\")) {\n        $pr2 = $p.Range\n        $tempFind = $pr2.Find\n        $tempFind.ClearFormatting()\n        $tempFind.Text = \"2.TEMP\"\n        $tempFound = $tempFind.Execute()\n        if ($tempFound) {\n            $pr2.Text = \"2.\"\n        }\n        break\n    }\n}\n"}
